# 006 E2E testcase for Hydroflask us
#
# Adds two new product rows (Wide Mouth Flex Sip(TM) Lid / 40 oz All
# Around(TM) Travel Tumbler) to the "E2E" sheet, just below the existing
# product rows, pushing the remaining test data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2E")

# Insert two new blank rows at row 16 (formats/styles are inherited from
# the row above, matching the surrounding AB/AC column styling).
$ws.Rows.Item(16).Resize(2).Insert()

# Row 16 - SKU-CFX454 / Wide Mouth Flex Sip(TM) Lid / Qty 1 / Laguna
$ws.Range("A16").Value = "SKU-CFX454"
$ws.Range("AB16").Value = "Wide Mouth Flex Sip™ Lid"
$ws.Range("AC16").Value = "'1"
$ws.Range("AE16").Value = "Laguna"

# Row 17 - SKU-TT40PS474 / 40 oz All Around(TM) Travel Tumbler / Qty 1 / Lupine
$ws.Range("A17").Value = "SKU-TT40PS474"
$ws.Range("AB17").Value = "40 oz All Around™ Travel Tumbler"
$ws.Range("AC17").Value = "'1"
$ws.Range("AE17").Value = "Lupine"

# Update the sheet's saved selection to match the edited area.
$ws.Range("AD14").Select()
